# Convert some images to code
#
# The deck had two "screenshot-style" code boxes (Content Placeholder 3,
# shape ids 7 and 10) duplicating the class Man{...} / class Woman{...}
# text that already exists as live, editable text boxes elsewhere on the
# slide. Remove those two redundant shapes along with the timeline
# animation effects (and build-list entries) that targeted them.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$targetIds = @(7, 10)

# 1) Remove the animation effects that target these shapes from the
#    slide's main timing sequence (this also drops the corresponding
#    <p:bldP> build-list entries). Walk backwards since deleting shifts
#    indices.
$tl = $s.TimeLine
$main = $tl.MainSequence
for ($i = $main.Count; $i -ge 1; $i--) {
    $eff = $main.Item($i)
    if ($targetIds -contains $eff.Shape.Id) {
        $eff.Delete()
    }
}

# 2) Remove the shapes themselves, again walking backwards.
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s.Shapes.Item($i)
    if ($targetIds -contains $shp.Id) {
        $shp.Delete()
    }
}
